# Apply the edits described in the diff:
# - absPath url shortened (folder removed)
# - workbookView window size/position changed
# - new log entries typed into rows 10 and 11 on sheet "DTT Test Hour Log"
# - selection moved from D11 to B5
# - SUMIF total in B30 recalculates from 14 to 16 as a natural consequence
#   of the new data (handled automatically by recalculation)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

# --- New log rows -----------------------------------------------------
# Columns A first (subject), matching the order the strings were
# originally authored in, so shared-string indices line up.
$ws.Range("A10").Value = "Createing Search function"
$ws.Range("A11").Value = "Added Pagination"

# Bonus flag for row 11
$ws.Range("E11").Value = "X"

# Description column for both new rows
$ws.Range("D10").Value = "I used the where clause methods to retrieve the facility information with multiple search queries"
$ws.Range("D11").Value = " created a pagination method to be used for api request"

# Amount of hours / date / bonus for row 10
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "09/26/2024"
$ws.Range("E10").Value = " "

# Amount of hours (text "0.5", matching the existing style used in B4)
# for row 11 - copy format+value from B4 so the cell keeps its numeric
# style id but is stored as the shared string "0.5" (text), not a number.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C11").Value = "09/26/2024"

# --- Selection ----------------------------------------------------------
$ws.Range("B5").Select()

# --- Window geometry ------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 5856
$win.Top = 0
$win.Width = 17280
$win.Height = 12336
